$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string for datelisted (affects C2 and C3 which both reference it)
$ws.Range("C2").Value = "29Jan2022"
$ws.Range("C3").Value = "29Jan2022"

# Update numeric values for row 2 and row 3 (duplicate rows, same new values)
$ws.Range("K2:K3").Value = 0
$ws.Range("L2:L3").Value = 1
$ws.Range("O2:O3").Value = 13.250951766967773
$ws.Range("P2:P3").Value = 38.024566650390625
$ws.Range("Q2:Q3").Value = 29.036596298217773
$ws.Range("R2:R3").Value = 8.9879703521728516
$ws.Range("S2:S3").Value = 3.2160265445709229
$ws.Range("T2:T3").Value = 35.548210144042969
$ws.Range("U2:U3").Value = 8.6075477600097656
$ws.Range("V2:V3").Value = 26.940662384033203
$ws.Range("X2:X3").Value = 0
$ws.Range("Y2:Y3").Value = 17.446767807006836
$ws.Range("Z2:Z3").Value = 21.880966186523438
$ws.Range("AA2:AA3").Value = 1.1773288249969482
$ws.Range("AB2:AB3").Value = 20.703638076782227
$ws.Range("AC2:AC3").Value = 8.8961515426635742
$ws.Range("AD2:AD3").Value = 2.1978754997253418
$ws.Range("AE2:AE3").Value = 2.1978754997253418
$ws.Range("AF2:AF3").Value = 0
$ws.Range("AG2:AG3").Value = 1
$ws.Range("AI2:AI3").Value = 7.0325741767883301
$ws.Range("AJ2:AJ3").Value = 29.231773376464844
$ws.Range("AK2:AK3").Value = 29.231773376464844
$ws.Range("AM2:AM3").Value = 4.6432280540466309
$ws.Range("AN2:AN3").Value = 8.4582910537719727
$ws.Range("AO2:AO3").Value = 8.4582910537719727
$ws.Range("AP2:AP3").Value = 0
$ws.Range("AQ2:AQ3").Value = 1
$ws.Range("AS2:AS3").Value = 2.3537311553955078
$ws.Range("AT2:AT3").Value = 32.880050659179688
$ws.Range("AU2:AU3").Value = 9.2514591217041016
$ws.Range("AV2:AV3").Value = 23.628591537475586
$ws.Range("AW2:AW3").Value = 8.6124334335327148
$ws.Range("AX2:AX3").Value = 13.234278678894043
$ws.Range("AY2:AY3").Value = 13.234278678894043
$ws.Range("BA2:BA3").Value = 1
$ws.Range("BC2:BC3").Value = 1.9467545747756958
$ws.Range("BD2:BD3").Value = 44.056102752685547
$ws.Range("BE2:BE3").Value = 39.308578491210938
$ws.Range("BF2:BF3").Value = 4.7475242614746094
$ws.Range("BG2:BG3").Value = 3.0712547302246094
$ws.Range("BH2:BH3").Value = 47.781539916992188
$ws.Range("BI2:BI3").Value = 43.099933624267578
$ws.Range("BJ2:BJ3").Value = 4.6816062927246094
$ws.Range("BM2:BM3").Value = 18.21574592590332
$ws.Range("BN2:BN3").Value = 42.9378662109375
$ws.Range("BO2:BO3").Value = 7.8077750205993652
$ws.Range("BP2:BP3").Value = 35.130092620849609
$ws.Range("BQ2:BQ3").Value = 14.384234428405762
$ws.Range("BR2:BR3").Value = 46.834102630615234
$ws.Range("BS2:BS3").Value = 3.9423618316650391
$ws.Range("BT2:BT3").Value = 42.891738891601563
$ws.Range("BU2:BU3").Value = 10.666536331176758
$ws.Range("BV2:BV3").Value = 7.4639077186584473
Write-Host "edit applied"
